$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.469.40"
$ws.Range("D3").Value = "2.285.25"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'309.73"
$ws.Range("E5").Value = "  -3.97%  "
$ws.Range("D6").Value = "'103.30"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  -1.21%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'38.53"
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("D11").Value = "'0.0899"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "'8.19"
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").Value = "'0.967"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "'15.16"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "2.630.72"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "2.291.67"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "42.433.41"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'13.48"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0000104"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").Value = "'73.00"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'265.28"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'3.39"
$ws.Range("E24").Value = "  -6.16%  "
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "'10.69"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.33"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'7.01"
$ws.Range("E29").Value = "  +13.94%  "
$ws.Range("D30").Value = "'21.81"
$ws.Range("E30").Value = "  -3.23%  "
$ws.Range("D31").Value = "'35.82"
$ws.Range("E31").Value = "  -5.36%  "
$ws.Range("D32").Value = "'164.49"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "'0.0849"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").Value = "'2.55"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("E36").Value = "  -3.53%  "
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("D38").Value = "'0.0345"
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("D39").Value = "'2.75"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "'3.62"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").Value = "'1.56"
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("D42").Value = "'102.85"
$ws.Range("E42").Value = "  +9.87%  "
$ws.Range("D43").Value = "'69.67"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").Value = "'12.03"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").Value = "1.712.59"
$ws.Range("E47").Value = "  +6.83%  "
$ws.Range("D48").Value = "'110.01"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "'77.56"
$ws.Range("E49").Value = "  -4.19%  "
$ws.Range("D50").Value = "'8.63"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").Value = "'5.14"
$ws.Range("E51").Value = "  -1.57%  "
